$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting rows 27:41 down to 28:42
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the new entry's data.
$ws.Cells.Item(27, 1).Value = 2
$ws.Cells.Item(27, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(27, 3).Value = "Coquimbo"
$ws.Cells.Item(27, 4).Value = 44489
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(27, 6).Value = 100112026
$ws.Cells.Item(27, 7).Value = "Haba"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 1200
$ws.Cells.Item(27, 11).Value = 5000
$ws.Cells.Item(27, 12).Value = 6000
$ws.Cells.Item(27, 13).Value = 5500
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 16).Value = 220
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
